# Merge the split "你好，" / "$\SunQuarTeX$" (and trailing "!") runs on
# slide 1 (Title) and slide 2 (Content) into a single run each, e.g.
# "你好，" + "$\SunQuarTeX$"      -> "你好，SunQuarTeX"
# "你好，" + "$\SunQuarTeX$" + "!" -> "你好，SunQuarTeX!"

$p = $ppt.ActivePresentation

# --- Slide 1: Title shape -------------------------------------------------
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange

# Drop the literal "$\SunQuarTeX$" run (chars 4-16) so only the "你好，"
# run (with its original, empty <a:rPr/>) remains.
$tr1.Characters(4, 13).Text = ""

# Re-fetch and append the replacement text onto the sole surviving run.
$tr1b = $sh1.TextFrame.TextRange
$tr1b.Text = $tr1b.Text + "SunQuarTeX"

# --- Slide 2: Content shape ------------------------------------------------
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(1)
$tr2 = $sh2.TextFrame.TextRange

# Drop the literal "$\SunQuarTeX$" run (chars 4-16), leaving "你好，" + "!".
$tr2.Characters(4, 13).Text = ""

# Drop the now-adjacent "!" run (char 4 of the remaining 4-char text too),
# leaving only the "你好，" run.
$tr2b = $sh2.TextFrame.TextRange
$tr2b.Characters(4, 1).Text = ""

# Re-fetch and append the replacement text onto the sole surviving run.
$tr2c = $sh2.TextFrame.TextRange
$tr2c.Text = $tr2c.Text + "SunQuarTeX!"
